$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6444
$ws.Range("C2").Value = 772
$ws.Range("D2").Value = 1490
$ws.Range("E2").Value = 81.6697703289882
$ws.Range("F2").Value = 85.22020725388602
$ws.Range("G2").Value = 92.78120805369127
$ws.Range("H2").Value = 3460
$ws.Range("I2").Value = 0.6574447062400244
$ws.Range("J2").Value = 2510
$ws.Range("K2").Value = 3.815169478644171
$ws.Range("L2").Value = 2766
$ws.Range("M2").Value = 2.000810161743005
$ws.Range("N2").Value = 1851
$ws.Range("O2").Value = 0.3517139165463252
$ws.Range("P2").Value = 230
$ws.Range("Q2").Value = 0.3495972032223743
$ws.Range("R2").Value = 405
$ws.Range("S2").Value = 0.2929602731402448

$ws.Range("B6").Value = 5936
$ws.Range("C6").Value = 383
$ws.Range("D6").Value = 767
$ws.Range("E6").Value = 155.2505053908356
$ws.Range("F6").Value = 96.18537859007833
$ws.Range("G6").Value = 96.77444589308996
$ws.Range("H6").Value = 1097
$ws.Range("I6").Value = 0.1190363804259484
$ws.Range("J6").Value = 785
$ws.Range("K6").Value = 2.13089388962784
$ws.Range("L6").Value = 990
$ws.Range("M6").Value = 1.333764449114865
$ws.Range("N6").Value = 140
$ws.Range("O6").Value = 0.01519151618927327
$ws.Range("P6").Value = 32
$ws.Range("Q6").Value = 0.08686446429056163
$ws.Range("R6").Value = 73
$ws.Range("S6").Value = 0.09834828766200523

$ws.Range("B11").Value = 7926
$ws.Range("C11").Value = 1132
$ws.Range("D11").Value = 2266
$ws.Range("E11").Value = 52.47451425687611
$ws.Range("F11").Value = 52.47614840989399
$ws.Range("G11").Value = 51.96028243601059
$ws.Range("H11").Value = 1862
$ws.Range("I11").Value = 0.4476897812763727
$ws.Range("J11").Value = 1333
$ws.Range("K11").Value = 2.243994411056681
$ws.Range("L11").Value = 1557
$ws.Range("M11").Value = 1.3223828370505
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0

$ws.Range("B12").Value = 4432
$ws.Range("C12").Value = 633
$ws.Range("D12").Value = 1267
$ws.Range("E12").Value = 29.50676895306859
$ws.Range("F12").Value = 26.64296998420221
$ws.Range("G12").Value = 28.53433307024467
$ws.Range("H12").Value = 2118
$ws.Range("I12").Value = 1.619587991496781
$ws.Range("J12").Value = 1334
$ws.Range("K12").Value = 7.909872517047139
$ws.Range("L12").Value = 1592
$ws.Range("M12").Value = 4.403507316128676
$ws.Range("N12").Value = 36
$ws.Range("O12").Value = 0.02752840778748069
$ws.Range("P12").Value = 8
$ws.Range("Q12").Value = 0.04743551734361103
$ws.Range("R12").Value = 6
$ws.Range("S12").Value = 0.01659613310098747

$ws.Range("B14").Value = 36000
$ws.Range("C14").Value = 1333
$ws.Range("D14").Value = 2667
$ws.Range("E14").Value = 26.3625
$ws.Range("F14").Value = 27.47861965491373
$ws.Range("G14").Value = 25.92688413948256
$ws.Range("H14").Value = 1095
$ws.Range("I14").Value = 0.115378536431168
$ws.Range("J14").Value = 575
$ws.Range("K14").Value = 1.569794425182232
$ws.Range("L14").Value = 521
$ws.Range("M14").Value = 0.7534672509291799
$ws.Range("N14").Value = 14830
$ws.Range("O14").Value = 1.562615246825773
$ws.Range("P14").Value = 588
$ws.Range("Q14").Value = 1.605285429577657
$ws.Range("R14").Value = 1140
$ws.Range("S14").Value = 1.648661547138704

$ws.Range("B16").Value = 4486
$ws.Range("C16").Value = 105
$ws.Range("D16").Value = 211
$ws.Range("E16").Value = 56.87895675434686
$ws.Range("F16").Value = 23.37142857142857
$ws.Range("G16").Value = 22.62085308056872
$ws.Range("H16").Value = 1597
$ws.Range("I16").Value = 0.6258842525640875
$ws.Range("J16").Value = 555
$ws.Range("K16").Value = 22.6161369193154
$ws.Range("L16").Value = 690
$ws.Range("M16").Value = 14.45631678189818
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0.0003919124937783891
$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 0.04074979625101875

$ws.Range("B17").Value = 2468
$ws.Range("C17").Value = 353
$ws.Range("D17").Value = 706
$ws.Range("E17").Value = 30.46839546191248
$ws.Range("F17").Value = 28.72804532577904
$ws.Range("G17").Value = 32.78895184135978
$ws.Range("H17").Value = 479
$ws.Range("I17").Value = 0.6370019681897974
$ws.Range("J17").Value = 527
$ws.Range("K17").Value = 5.196726161128094
$ws.Range("L17").Value = 572
$ws.Range("M17").Value = 2.470949069074258
$ws.Range("N17").Value = 624
$ws.Range("O17").Value = 0.8298313740092558
$ws.Range("P17").Value = 76
$ws.Range("Q17").Value = 0.749432994773691
$ws.Range("R17").Value = 147
$ws.Range("S17").Value = 0.6350166313879648

$ws.Range("B20").Value = 595
$ws.Range("C20").Value = 85
$ws.Range("D20").Value = 171
$ws.Range("E20").Value = 54.32100840336135
$ws.Range("F20").Value = 52.83529411764706
$ws.Range("G20").Value = 57.8187134502924
$ws.Range("H20").Value = 1195
$ws.Range("I20").Value = 3.697286593855388
$ws.Range("J20").Value = 588
$ws.Range("K20").Value = 13.09285237140949
$ws.Range("L20").Value = 807
$ws.Range("M20").Value = 8.162233235561848
$ws.Range("N20").Value = 235
$ws.Range("O20").Value = 0.727081464063612
$ws.Range("P20").Value = 30
$ws.Range("Q20").Value = 0.6680026720106881
$ws.Range("R20").Value = 108
$ws.Range("S20").Value = 1.092343481339132

